$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------------
# Rows 19-26: "Poursuivre le systeme de vie" section - fill in the
# task/duration table that up to now had only the A column (week number)
# filled in, plus remove the now-unused B26 label cell.
#
# Shared strings must be appended in the exact order the new text values are
# assigned below, matching sharedStrings.xml indices 30..37 from the target
# workbook.
# ---------------------------------------------------------------------------

# Row 19 - section header
$rng = $ws.Range("B19")
$rng.Value = "Poursuivre le système de vie"
$rng.Font.Size = 12
$rng.Font.Color = 0
$rng.Interior.Color = 14277081
$rng.Interior.PatternColor = 0

# Row 20
$rng = $ws.Range("B20")
$rng.Value = "          Faire réagir dégats/baisse de la barre de vie"
$rng.Font.Size = 12
$rng.Font.Color = 0
$rng.Interior.Color = 16777215
$rng.Interior.PatternColor = 0

# Row 23
$rng = $ws.Range("B23")
$rng.Value = "          Créer l'évènement ""mort"" et afficher écran Game Over"
$rng.Font.Size = 12
$rng.Font.Color = 0
$rng.Interior.Color = 16777215
$rng.Interior.PatternColor = 0

# Row 25
$rng = $ws.Range("B25")
$rng.Value = "Afficher le background"
$rng.Font.Size = 12
$rng.Font.Color = 0

# Row 22
$rng = $ws.Range("B22")
$rng.Value = "                     Ajout du statut Current Health et Max Health dans le script ""Health Bar"""
$rng.Borders.LineStyle = -4142

# Row 24
$rng = $ws.Range("B24")
$rng.Value = "                     Ajout de la condition if(Current Health == 0)"
$rng.Font.Size = 12
$rng.Font.Color = 0
$rng.Interior.Color = 16777215
$rng.Interior.PatternColor = 0

# Row 21
$rng = $ws.Range("B21")
$rng.Value = "                     Création du scripth ""Health Bar"""
$rng.Borders.LineStyle = -4142

# Row 26 - the label cell in column B is removed entirely
$rng = $ws.Range("B26")
$rng.Style = "Normal"
$rng.ClearContents()

# C column values (durations). C19's "6,5" is a shared string and must be
# set last so it lands on shared-string index 37.
$ws.Range("C20").Value = 3
$ws.Range("C21").Value = 2
$ws.Range("C22").Value = 1
$ws.Range("C23").Value = "0,5"
$ws.Range("C24").Value = "0,5"
$ws.Range("C25").Value = "0,5"
$ws.Range("C19").Value = "6,5"

# ---------------------------------------------------------------------------
# Update the active selection shown when the workbook is reopened.
# ---------------------------------------------------------------------------
$ws.Range("C26").Select()
